# Atualizado por script em 01-12-2023 20:45
#
# Applies:
#  1. Swap the F:V (match-detail) columns between row 8 and row 9
#     (same kickoff date, order of the two fixtures corrected).
#  2. Swap the F:V (match-detail) columns between row 12 and row 13.
#  3. Swap the F:V (match-detail) columns between row 41 and row 42.
#  4. Append three new fixtures as rows 55, 56 and 57.
#
# NOTE: this interpreter only binds function parameters *positionally* -
# named parameters (-Foo bar) are silently ignored - so every helper below
# takes its arguments by position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowDetails {
    param(
        $Row,
        $F, $G, $H, $I,
        $J, $K, $L, $M, $N, $O, $P, $Q, $R, $S, $T, $U,
        $V
    )

    $ws.Cells.Item($Row, 6).Value2  = $F
    $ws.Cells.Item($Row, 7).Value2  = $G
    $ws.Cells.Item($Row, 8).Value2  = $H
    $ws.Cells.Item($Row, 9).Value2  = $I
    $ws.Cells.Item($Row, 10).Value2 = $J
    $ws.Cells.Item($Row, 11).Value2 = $K
    $ws.Cells.Item($Row, 12).Value2 = $L
    $ws.Cells.Item($Row, 13).Value2 = $M
    $ws.Cells.Item($Row, 14).Value2 = $N
    $ws.Cells.Item($Row, 15).Value2 = $O
    $ws.Cells.Item($Row, 16).Value2 = $P
    $ws.Cells.Item($Row, 17).Value2 = $Q
    $ws.Cells.Item($Row, 18).Value2 = $R
    $ws.Cells.Item($Row, 19).Value2 = $S
    $ws.Cells.Item($Row, 20).Value2 = $T
    $ws.Cells.Item($Row, 21).Value2 = $U
    $ws.Cells.Item($Row, 22).Value2 = $V
}

# --- capture the current ("before") F:V values of a row so the writes
#     below don't clobber data still to be read during a swap.
function Get-RowDetails {
    param($Row)

    [PSCustomObject]@{
        F = $ws.Cells.Item($Row, 6).Value2
        G = $ws.Cells.Item($Row, 7).Value2
        H = $ws.Cells.Item($Row, 8).Value2
        I = $ws.Cells.Item($Row, 9).Value2
        J = $ws.Cells.Item($Row, 10).Value2
        K = $ws.Cells.Item($Row, 11).Value2
        L = $ws.Cells.Item($Row, 12).Value2
        M = $ws.Cells.Item($Row, 13).Value2
        N = $ws.Cells.Item($Row, 14).Value2
        O = $ws.Cells.Item($Row, 15).Value2
        P = $ws.Cells.Item($Row, 16).Value2
        Q = $ws.Cells.Item($Row, 17).Value2
        R = $ws.Cells.Item($Row, 18).Value2
        S = $ws.Cells.Item($Row, 19).Value2
        T = $ws.Cells.Item($Row, 20).Value2
        U = $ws.Cells.Item($Row, 21).Value2
        V = $ws.Cells.Item($Row, 22).Value2
    }
}

function Swap-RowDetails {
    param($RowA, $RowB)

    $a = Get-RowDetails $RowA
    $b = Get-RowDetails $RowB

    Set-RowDetails $RowA $b.F $b.G $b.H $b.I $b.J $b.K $b.L $b.M $b.N $b.O $b.P $b.Q $b.R $b.S $b.T $b.U $b.V
    Set-RowDetails $RowB $a.F $a.G $a.H $a.I $a.J $a.K $a.L $a.M $a.N $a.O $a.P $a.Q $a.R $a.S $a.T $a.U $a.V
}

# 1) rows 8 / 9
Swap-RowDetails 8 9

# 2) rows 12 / 13
Swap-RowDetails 12 13

# 3) rows 41 / 42
Swap-RowDetails 41 42

# 4) append new fixtures as rows 55-57, copying the style (bold/border on
#    col A, date-time number format on col E) from the last existing row.
function Add-FixtureRow {
    param(
        $Row, $Indice, $DataPartida,
        $F, $G, $H, $I,
        $J, $K, $L, $M, $N, $O, $P, $Q, $R, $S, $T, $U,
        $V
    )

    $ws.Cells.Item(54, 1).Copy()
    $ws.Cells.Item($Row, 1).PasteSpecial(-4122)
    $ws.Cells.Item(54, 5).Copy()
    $ws.Cells.Item($Row, 5).PasteSpecial(-4122)
    $excel.CutCopyMode = 0

    $ws.Cells.Item($Row, 1).Value2 = $Indice
    $ws.Cells.Item($Row, 2).Value2 = "algeria"
    $ws.Cells.Item($Row, 3).Value2 = "ligue-1"
    $ws.Cells.Item($Row, 4).Value2 = "2023-2024"
    $ws.Cells.Item($Row, 5).Value2 = $DataPartida

    Set-RowDetails $Row $F $G $H $I $J $K $L $M $N $O $P $Q $R $S $T $U $V
}

Add-FixtureRow 55 54 45261.63541666666 `
    "El Bayadh" 1 "Paradou" 0 `
    1.95 "30/11/2023 07:50" 2.06 "01/12/2023 15:11" `
    3.08 "30/11/2023 07:50" 2.96 "01/12/2023 15:11" `
    4.28 "30/11/2023 07:50" 4.27 "01/12/2023 15:11" `
    "https://www.betexplorer.com/football/algeria/ligue-1/el-bayadh-paradou/U99Qf5Er/"

Add-FixtureRow 56 55 45261.69791666666 `
    "MC Alger" 4 "Magra" 0 `
    1.35 "30/11/2023 07:50" 1.22 "01/12/2023 15:57" `
    4.62 "30/11/2023 07:50" 5.76 "01/12/2023 15:57" `
    9.279999999999999 "30/11/2023 07:50" 15.13 "01/12/2023 15:57" `
    "https://www.betexplorer.com/football/algeria/ligue-1/mc-alger-magra/bqEUgPTl/"

Add-FixtureRow 57 56 45261.75 `
    "Constantine" 2 "Khenchela" 0 `
    1.76 "30/11/2023 06:12" 1.61 "01/12/2023 17:55" `
    3.13 "30/11/2023 06:12" 3.49 "01/12/2023 17:55" `
    4.92 "30/11/2023 06:12" 6.62 "01/12/2023 17:55" `
    "https://www.betexplorer.com/football/algeria/ligue-1/constantine-khenchela/S0Nbao6R/"

Write-Host "Edit applied."
